$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-parsed as numbers by Excel, so they stay text like the source data.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the refreshed coin data (rank order, price, 1h volume change).
$ws.Range("D2").Value = '69.542.95'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '3.744.03'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").Value = '608.73'
$ws.Range("E5").Value = '  +3.31%  '
$ws.Range("D6").Value = '187.27'
$ws.Range("E6").Value = '  +14.23%  '
$ws.Range("D7").Value = '3.735.51'
$ws.Range("E7").Value = '  -1.91%  '
$ws.Range("D8").Value = '0.638'
$ws.Range("E8").Value = '  -4.06%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '0.727'
$ws.Range("E10").Value = '  -1.49%  '
$ws.Range("D11").Value = '0.166'
$ws.Range("E11").Value = '  -3.75%  '
$ws.Range("D12").Value = '56.88'
$ws.Range("E12").Value = '  +8.65%  '
$ws.Range("D13").Value = '0.0000298'
$ws.Range("E13").Value = '  -5.59%  '
$ws.Range("D14").Value = '10.61'
$ws.Range("E14").Value = '  -5.05%  '
$ws.Range("D15").Value = '4.308.55'
$ws.Range("E15").Value = '  -3.08%  '
$ws.Range("D16").Value = '3.715.57'
$ws.Range("E16").Value = '  -3.21%  '
$ws.Range("D17").Value = '19.62'
$ws.Range("E17").Value = '  -5.20%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '13.02'
$ws.Range("E18").Value = '  -4.56%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.127'
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("D20").Value = '1.14'
$ws.Range("E20").Value = '  -4.71%  '
$ws.Range("D21").Value = '69.005.66'
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("D22").Value = '415.82'
$ws.Range("E22").Value = '  -3.74%  '
$ws.Range("D23").Value = '4.71'
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("D24").Value = '90.08'
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("D25").Value = '3.07'
$ws.Range("E25").Value = '  -4.68%  '
$ws.Range("D26").Value = '13.03'
$ws.Range("E26").Value = '  -4.93%  '
$ws.Range("D27").Value = '11.15'
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").Value = '3.96'
$ws.Range("E28").Value = '  +1.27%  '
$ws.Range("D29").Value = '6.07'
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("D30").Value = '9.63'
$ws.Range("E30").Value = '  -6.70%  '
$ws.Range("D31").Value = '33.25'
$ws.Range("E31").Value = '  -4.15%  '
$ws.Range("D32").Value = '7.46'
$ws.Range("E32").Value = '  -6.77%  '
$ws.Range("D33").Value = '12.64'
$ws.Range("E33").Value = '  -4.93%  '
$ws.Range("D34").Value = '0.119'
$ws.Range("E34").Value = '  -4.01%  '
$ws.Range("D35").Value = '44.48'
$ws.Range("E35").Value = '  -6.52%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = '611.61'
$ws.Range("E36").Value = '  -3.16%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '64.98'
$ws.Range("E37").Value = '  -6.16%  '
$ws.Range("D38").Value = '0.0₃0910'
$ws.Range("E38").Value = '  -5.98%  '
$ws.Range("D39").Value = '0.407'
$ws.Range("E39").Value = '  -2.68%  '
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("D41").Value = '0.996'
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("D42").Value = '0.138'
$ws.Range("E42").Value = '  -3.80%  '
$ws.Range("D43").Value = '3.08'
$ws.Range("E43").Value = '  -4.29%  '
$ws.Range("D44").Value = '2.79'
$ws.Range("E44").Value = '  +3.83%  '
$ws.Range("D45").Value = '0.0445'
$ws.Range("E45").Value = '  -4.08%  '
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  -6.17%  '
$ws.Range("D47").Value = '9.39'
$ws.Range("E47").Value = '  -4.26%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.136'
$ws.Range("E48").Value = '  -4.41%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = '2.75'
$ws.Range("E49").Value = '  -3.00%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.804.52'
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '3.16'
$ws.Range("E51").Value = '  -1.91%  '
